# Update the projection results on the active sheet with refreshed
# participation / contribution figures for rows 2-6 (years 1-5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year 1 (row 2)
$ws.Range("C2").Value = 9448
$ws.Range("D2").Value = 7574
$ws.Range("E2").Value = 0.8016511430990686
$ws.Range("F2").Value = 0.7991137370753324
$ws.Range("G2").Value = 0.1028346976498548
$ws.Range("H2").Value = 0.08217661953998734
$ws.Range("I2").Value = 39239770.164845
$ws.Range("J2").Value = 13420220.8417585
$ws.Range("L2").Value = 13420220.8417585
$ws.Range("M2").Value = 52659991.0066035
$ws.Range("N2").Value = 801375021.0472001
$ws.Range("O2").Value = 783675214.0432
$ws.Range("P2").Value = 0.01674649257749708
$ws.Range("Q2").Value = 0.01712472284598593

# Year 2 (row 3)
$ws.Range("C3").Value = 9640
$ws.Range("D3").Value = 7551
$ws.Range("E3").Value = 0.783298755186722
$ws.Range("F3").Value = 0.7811918063314711
$ws.Range("G3").Value = 0.1028830618461131
$ws.Range("H3").Value = 0.08037140492447756
$ws.Range("I3").Value = 40538728.96122567
$ws.Range("J3").Value = 13799638.26741243
$ws.Range("L3").Value = 13799638.26741243
$ws.Range("M3").Value = 54338367.2286381
$ws.Range("N3").Value = 836150914.9537281
$ws.Range("O3").Value = 818670738.9296581
$ws.Range("P3").Value = 0.01650376507472469
$ws.Range("Q3").Value = 0.01685615182173761

# Year 3 (row 4)
$ws.Range("C4").Value = 9834
$ws.Range("D4").Value = 7543
$ws.Range("E4").Value = 0.7670327435428107
$ws.Range("F4").Value = 0.7651653479407587
$ws.Range("G4").Value = 0.1029497547394936
$ws.Range("H4").Value = 0.07877358490566039
$ws.Range("I4").Value = 42014448.3095379
$ws.Range("J4").Value = 14217426.49631654
$ws.Range("L4").Value = 14217426.49631654
$ws.Range("M4").Value = 56231874.80585443
$ws.Range("N4").Value = 875322312.16536
$ws.Range("O4").Value = 857873364.1594061
$ws.Range("P4").Value = 0.0162425043880644
$ws.Range("Q4").Value = 0.01657287321217578

# Year 4 (row 5)
$ws.Range("C5").Value = 10026
$ws.Range("D5").Value = 7534
$ws.Range("E5").Value = 0.7514462397765809
$ws.Range("F5").Value = 0.7493534911478019
$ws.Range("G5").Value = 0.1030222989116007
$ws.Range("H5").Value = 0.07720011935548041
$ws.Range("I5").Value = 43535019.83059579
$ws.Range("J5").Value = 14644484.75463068
$ws.Range("L5").Value = 14644484.75463068
$ws.Range("M5").Value = 58179504.58522647
$ws.Range("N5").Value = 914181328.8136762
$ws.Range("O5").Value = 896695223.350703
$ws.Range("P5").Value = 0.01601923414213095
$ws.Range("Q5").Value = 0.01633161900863961

# Year 5 (row 6)
$ws.Range("C6").Value = 10237
$ws.Range("D6").Value = 7523
$ws.Range("E6").Value = 0.7348832665820064
$ws.Range("F6").Value = 0.7336649112541447
$ws.Range("G6").Value = 0.1030931809118703
$ws.Range("H6").Value = 0.07563584942461479
$ws.Range("I6").Value = 45146428.42506469
$ws.Range("J6").Value = 15078985.98419153
$ws.Range("L6").Value = 15078985.98419153
$ws.Range("M6").Value = 60225414.4092562
$ws.Range("N6").Value = 955175919.0364679
$ws.Range("O6").Value = 937584093.1632864
$ws.Range("P6").Value = 0.01578660609388313
$ws.Range("Q6").Value = 0.01608280909855989
